# Add 9 new user rows (110021 - 110029) to the master-user_detail sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(110021, 7316931025, "Magdalena Weber",   "magdalena.weber@xyz.com",   932122450),
    @(110022, 9137847236, "Adrienne Hoffman",  "adrienne.hoffman@xyz.com",  848488000),
    @(110023, 8428758532, "Adrienne Mcgee",    "adrienne.mcgee@xyz.com",    894773246),
    @(110024, 9804209494, "Amare Coleman",     "amare.coleman@xyz.com",     956554588),
    @(110025, 7105248214, "Dawson Ibarra",     "dawson.ibarra@xyz.com",     765455583),
    @(110026, 9316557128, "Elvis Mcmillan",    "elvis.mcmillan@xyz.com",    884282274),
    @(110027, 8103486949, "Steve George",      "steve.george@xyz.com",      971073663),
    @(110028, 9601932866, "Colton Elliott",    "colton.elliott@xyz.com",    809908673),
    @(110029, 9317596765, "Carolyn Rodriguez", "carolyn.rodriguez@xyz.com", 818876429)
)

$startRow = 22

# The shared-strings table is built in first-seen order as cells are
# populated. The source workbook lists all nine new names before any of
# the new e-mail addresses, so fill column by column (not row by row) to
# reproduce that same ordering: id, uin, name(all), email(all), mobile...
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $data[$i][2]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $data[$i][3]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $data[$i][4]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = "ACT"
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = "eng"
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 8).Value = "PWD"
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 9).Value = $true
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value = "superadmin"
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 11).Value = "now()"
}

# Column D (e-mail) carries the same "alignment: left" style as the rows
# above it; column I (is_active) carries the boolean's style too.
$ws.Range("D22:D30").Style = $ws.Range("D21").Style
$ws.Range("I22:I30").Style = $ws.Range("I21").Style

# Match the updated view state recorded in the author's session.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A16").Select()
$ws.Range("A22:K30").Select()
